$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: columns A-G for rows 3-6 (values taken from the source dataset)
$rowsData = @(
    ,@(1, ' 24/08/2021', -0.0001925492790740479, -0.003411609458509669, -0.0007445663914916307, 0.004273972602739726, '{''ACE.NZ'': 0.0, ''AFC.NZ'': 0.0, ''AFI.NZ'': 0.0, ''AFT.NZ'': 0.0, ''AGL.NZ'': -0.006134969325153228, ''AIA.NZ'': -0.005486968449931462, ''AIR.NZ'': 0.01388888888888884, ''ALF.NZ'': -0.016949152542372947, ''AMP.NZ'': 0.0, ''ANZ.NZ'': -0.0003374957813027102, ''APL.NZ'': 0.015151515151515138, ''ARB.NZ'': 0.015625, ''ARG.NZ'': 0.0030674846625766694, ''ARV.NZ'': 0.009708737864077666, ''ATM.NZ'': 0.0, ''BFG.NZ'': -0.025000000000000022, ''BGI.NZ'': 0.0, ''BGP.NZ'': 0.0015748031496063408, ''BIT.NZ'': 0.0, ''BLT.NZ'': 0.014925373134328401, ''BRM.NZ'': 0.009803921568627416, ''CAV.NZ'': 0.0, ''CBD.NZ'': 0.012195121951219523, ''CDI.NZ'': 0.0, ''CEN.NZ'': 0.008536585365853666, ''CGF.NZ'': 0.0, ''CMO.NZ'': 0.0, ''CNU.NZ'': -0.024531024531024515, ''CRP.NZ'': 0.02400000000000002, ''CVT.NZ'': 0.015151515151515138, ''DGL.NZ'': 0.006917755572636519, ''EBO.NZ'': 0.01522388059701485, ''ENS.NZ'': -0.017064846416382395, ''ERD.NZ'': 0.004807692307692291, ''EVO.NZ'': -0.03125, ''FBU.NZ'': -0.005181347150259086, ''FCT.NZ'': 0.0, ''FPH.NZ'': -0.002089552238806025, ''FSF.NZ'': 0.010810810810810922, ''FWL.NZ'': -0.013333333333333308, ''GEN.NZ'': -0.03508771929824561, ''GEO.NZ'': -0.01666666666666672, ''GFL.NZ'': 0.0, ''GMT.NZ'': 0.005769230769230749, ''GNE.NZ'': 0.010494752623688264, ''GSH.NZ'': 0.0, ''GTK.NZ'': 0.0, ''GWC.NZ'': 0.0, ''GXH.NZ'': 0.009009009009008917, ''HFL.NZ'': -0.0064102564102563875, ''HGH.NZ'': 0.033333333333333215, ''HLG.NZ'': 0.005952380952380931, ''IFT.NZ'': -0.0013368983957220415, ''IKE.NZ'': 0.03846153846153855, ''IPL.NZ'': 0.0, ''JLG.NZ'': 0.0, ''JPG.NZ'': 0.0, ''KFL.NZ'': 0.004975124378109541, ''KMD.NZ'': 0.0, ''KPG.NZ'': 0.004385964912280826, ''LIC.NZ'': 0.0, ''MCK.NZ'': -0.032653061224489854, ''MCY.NZ'': 0.01319648093841641, ''MEL.NZ'': -0.0018939393939395588, ''MFT.NZ'': 0.004888345739362254, ''MHJ.NZ'': 0.011363636363636465, ''MHM.NZ'': 0.015873015873015817, ''MLN.NZ'': 0.013513513513513598, ''MMH.NZ'': 0.0, ''MPG.NZ'': 0.0, ''MWE.NZ'': 0.0, ''NPH.NZ'': 0.012539184952978122, ''NTL.NZ'': 0.0, ''NWF.NZ'': 0.0185185185185186, ''NZK.NZ'': 0.014184397163120588, ''NZM.NZ'': 0.010752688172043001, ''NZO.NZ'': 0.012499999999999956, ''NZR.NZ'': 0.034482758620689724, ''NZX.NZ'': 0.010810810810810922, ''OCA.NZ'': 0.013698630136986356, ''PCT.NZ'': 0.0, ''PEB.NZ'': 0.0, ''PFI.NZ'': 0.010101010101009944, ''PGW.NZ'': 0.0437317784256559, ''PHL.NZ'': 0.0, ''POT.NZ'': -0.005540166204986208, ''PPH.NZ'': 0.03571428571428581, ''PX1.NZ'': -0.01754385964912264, ''PYS.NZ'': 0.03125, ''QEX.NZ'': 0.0, ''RAK.NZ'': 0.010526315789473717, ''RBD.NZ'': 0.013548387096774306, ''RYM.NZ'': 0.0026791694574681557, ''SAN.NZ'': 0.010638297872340496, ''SCL.NZ'': 0.02277904328018243, ''SCT.NZ'': -0.003424657534246478, ''SCY.NZ'': 0.0, ''SDL.NZ'': 0.0, ''SKC.NZ'': 0.013201320132013139, ''SKL.NZ'': -0.0018382352941177516, ''SKO.NZ'': -0.018867924528301883, ''SKT.NZ'': 0.018867924528301883, ''SML.NZ'': 0.020588235294117796, ''SNC.NZ'': 0.0, ''SPG.NZ'': -0.015384615384615441, ''SPK.NZ'': 0.003121748178980255, ''SPN.NZ'': 0.0023529411764704466, ''SPY.NZ'': -0.0326797385620915, ''STU.NZ'': 0.01785714285714257, ''SUM.NZ'': 0.05681818181818188, ''TCL.NZ'': 0.003821656050955413, ''TEM.NZ'': 0.008450704225352101, ''TGG.NZ'': 0.006825938566553003, ''THL.NZ'': 0.013215859030836885, ''TPW.NZ'': -0.0036363636363635488, ''TRA.NZ'': -0.006976744186046546, ''TRU.NZ'': -0.04761904761904767, ''TWR.NZ'': 0.014084507042253502, ''VCT.NZ'': 0.007352941176470562, ''VGL.NZ'': 0.01388888888888884, ''VHP.NZ'': 0.004573170731707377, ''VTL.NZ'': 0.0, ''WBC.NZ'': 0.005192878338278861, ''WDT.NZ'': 0.0, ''WHS.NZ'': 0.008333333333333304, ''ZEL.NZ'': -0.002873563218390718}')
    ,@(2, ' 28/08/2021', -0.004358792098968958, -0.001440108165108698, 0.001132076890713406, 0.004986301369863014, '{''ACE.NZ'': 0.0, ''AFC.NZ'': 0.0, ''AFI.NZ'': -0.0034207525655644, ''AFT.NZ'': -0.011494252873563204, ''AGL.NZ'': 0.024844720496894235, ''AIA.NZ'': -0.0070224719101122934, ''AIR.NZ'': -0.00666666666666671, ''ALF.NZ'': -0.016129032258064502, ''AMP.NZ'': 0.0, ''ANZ.NZ'': -0.003712453594330034, ''APL.NZ'': -0.014925373134328401, ''ARB.NZ'': -0.032258064516129115, ''ARG.NZ'': -0.006134969325153228, ''ARV.NZ'': 0.004878048780487809, ''ATM.NZ'': -0.039432176656151396, ''BFG.NZ'': -0.012658227848101222, ''BGI.NZ'': 0.0, ''BGP.NZ'': -0.0015384615384614886, ''BIT.NZ'': -0.029787234042553346, ''BLT.NZ'': 0.0, ''BRM.NZ'': 0.0, ''CAV.NZ'': 0.01538461538461533, ''CDI.NZ'': 0.01754385964912286, ''CEN.NZ'': -0.008663366336633671, ''CGF.NZ'': 0.026315789473684292, ''CMO.NZ'': 0.0, ''CNU.NZ'': 0.013138686131386912, ''CRP.NZ'': 0.0, ''CVT.NZ'': 0.0, ''DGL.NZ'': 0.015602836879432758, ''EBO.NZ'': 0.0005712653527565337, ''ENS.NZ'': -0.018181818181818077, ''ERD.NZ'': 0.0031847133757960666, ''EVO.NZ'': 0.0, ''FBU.NZ'': -0.02785145888594165, ''FCT.NZ'': -0.02506963788300831, ''FPH.NZ'': 0.008442776735459567, ''FRE.NZ'': -0.011032308904649235, ''FSF.NZ'': 0.0, ''FWL.NZ'': 0.013605442176870763, ''GEN.NZ'': 0.0, ''GEO.NZ'': -0.025423728813559254, ''GFL.NZ'': 0.01388888888888884, ''GMT.NZ'': 0.011450381679389166, ''GNE.NZ'': 0.0, ''GSH.NZ'': -0.014084507042253391, ''GTK.NZ'': -0.00990099009900991, ''GWC.NZ'': 0.0, ''GXH.NZ'': 0.017699115044247815, ''HFL.NZ'': 0.0, ''HGH.NZ'': 0.008928571428571175, ''HLG.NZ'': -0.0098314606741573, ''IFT.NZ'': 0.0054421768707482165, ''IKE.NZ'': 0.0, ''IPL.NZ'': -0.009803921568627416, ''JLG.NZ'': -0.010869565217391353, ''JPG.NZ'': -0.004319654427645703, ''KFL.NZ'': 0.0, ''KMD.NZ'': -0.014705882352941235, ''KPG.NZ'': 0.008658008658008587, ''LIC.NZ'': 0.0, ''MCK.NZ'': -0.012295081967213073, ''MCY.NZ'': 0.0, ''MEL.NZ'': 0.009960159362550014, ''MFT.NZ'': 0.005440104450005423, ''MHJ.NZ'': 0.0, ''MHM.NZ'': 0.0, ''MLN.NZ'': -0.019607843137254943, ''MMH.NZ'': 0.027732463295269127, ''MPG.NZ'': -0.011235955056179803, ''MWE.NZ'': 0.04081632653061229, ''NTL.NZ'': 0.0, ''NWF.NZ'': -0.017857142857142905, ''NZK.NZ'': -0.006896551724137945, ''NZM.NZ'': 0.010638297872340496, ''NZR.NZ'': 0.0, ''NZX.NZ'': -0.021505376344086002, ''OCA.NZ'': 0.006666666666666599, ''PCT.NZ'': -0.002941176470588225, ''PEB.NZ'': -0.03597122302158262, ''PFI.NZ'': -0.009900990099009799, ''PGW.NZ'': -0.010752688172043001, ''PHL.NZ'': 0.0, ''POT.NZ'': 0.03216783216783203, ''PPH.NZ'': 0.005714285714285783, ''PX1.NZ'': -0.05263157894736825, ''PYS.NZ'': -0.015384615384615441, ''QEX.NZ'': 0.0, ''RAK.NZ'': -0.020618556701030966, ''RBD.NZ'': 0.024983984625240208, ''RYM.NZ'': 0.011235955056179803, ''SAN.NZ'': -0.008620689655172376, ''SCL.NZ'': 0.030303030303030276, ''SCT.NZ'': 0.03191489361702149, ''SCY.NZ'': 0.0, ''SDL.NZ'': 0.03333333333333344, ''SEK.NZ'': 0.003824091778202643, ''SKC.NZ'': 0.01834862385321112, ''SKL.NZ'': 0.03612167300380231, ''SKO.NZ'': -0.0024937655860348684, ''SKT.NZ'': 0.01840490797546024, ''SML.NZ'': -0.015243902439024293, ''SPG.NZ'': -0.015151515151515138, ''SPK.NZ'': 0.0, ''SPN.NZ'': 0.008235294117647118, ''SPY.NZ'': 0.01333333333333342, ''STU.NZ'': 0.0, ''SUM.NZ'': 0.009210526315789558, ''TEM.NZ'': 0.02556818181818188, ''TGG.NZ'': 0.0, ''THL.NZ'': 0.04721030042918439, ''TLL.NZ'': 0.0, ''TLT.NZ'': 0.0, ''TPW.NZ'': -0.018633540372670843, ''TRA.NZ'': 0.0, ''TRU.NZ'': 0.03333333333333344, ''TWR.NZ'': -0.0071428571428571175, ''VCT.NZ'': -0.025821596244131273, ''VGL.NZ'': 0.08370044052863435, ''VHP.NZ'': -0.006069802731411222, ''VTL.NZ'': 0.013513513513513598, ''WBC.NZ'': -0.0007382798080471931, ''WDT.NZ'': -0.03149606299212604, ''WHS.NZ'': -0.0026385224274406704, ''ZEL.NZ'': 0.0}')
    ,@(3, ' 28/08/2021', -0.004358792098968958, -0.0004019093367033248, 0.004693009067110755, 0.004986301369863014, '{''ACE.NZ'': 0.0, ''AFC.NZ'': 0.0, ''AFI.NZ'': -0.0034207525655644, ''AFT.NZ'': -0.011494252873563204, ''AGL.NZ'': 0.024844720496894235, ''AIA.NZ'': -0.0070224719101122934, ''AIR.NZ'': -0.00666666666666671, ''ALF.NZ'': -0.016129032258064502, ''AMP.NZ'': 0.0, ''ANZ.NZ'': -0.003712453594330034, ''APL.NZ'': -0.014925373134328401, ''ARB.NZ'': -0.032258064516129115, ''ARG.NZ'': -0.006134969325153228, ''ARV.NZ'': 0.004878048780487809, ''ATM.NZ'': -0.039432176656151396, ''BFG.NZ'': -0.012658227848101222, ''BGI.NZ'': 0.0, ''BGP.NZ'': -0.0015384615384614886, ''BIT.NZ'': -0.029787234042553346, ''BLT.NZ'': 0.0, ''BRM.NZ'': 0.0, ''CAV.NZ'': 0.01538461538461533, ''CBD.NZ'': 0.012195121951219523, ''CDI.NZ'': 0.01754385964912286, ''CEN.NZ'': -0.008663366336633671, ''CGF.NZ'': 0.026315789473684292, ''CMO.NZ'': 0.0, ''CNU.NZ'': 0.013138686131386912, ''CRP.NZ'': 0.0, ''CVT.NZ'': 0.0, ''DGL.NZ'': 0.015602836879432758, ''EBO.NZ'': 0.0005712653527565337, ''ENS.NZ'': -0.018181818181818077, ''ERD.NZ'': 0.0031847133757960666, ''EVO.NZ'': 0.0, ''FBU.NZ'': -0.02785145888594165, ''FCT.NZ'': -0.02506963788300831, ''FPH.NZ'': 0.008442776735459567, ''FRE.NZ'': -0.011032308904649235, ''FSF.NZ'': 0.0, ''FWL.NZ'': 0.013605442176870763, ''GEN.NZ'': 0.0, ''GEO.NZ'': -0.025423728813559254, ''GFL.NZ'': 0.01388888888888884, ''GMT.NZ'': 0.011450381679389166, ''GNE.NZ'': 0.0, ''GSH.NZ'': -0.014084507042253391, ''GTK.NZ'': -0.00990099009900991, ''GWC.NZ'': 0.0, ''GXH.NZ'': 0.017699115044247815, ''HFL.NZ'': 0.0, ''HGH.NZ'': 0.008928571428571175, ''HLG.NZ'': -0.0098314606741573, ''IFT.NZ'': 0.0054421768707482165, ''IKE.NZ'': 0.0, ''IPL.NZ'': -0.009803921568627416, ''JLG.NZ'': -0.010869565217391353, ''JPG.NZ'': -0.004319654427645703, ''KFL.NZ'': 0.0, ''KMD.NZ'': -0.014705882352941235, ''KPG.NZ'': 0.008658008658008587, ''LIC.NZ'': 0.0, ''MCK.NZ'': -0.012295081967213073, ''MCY.NZ'': 0.0, ''MEL.NZ'': 0.009960159362550014, ''MFT.NZ'': 0.005440104450005423, ''MHJ.NZ'': 0.0, ''MHM.NZ'': 0.0, ''MLN.NZ'': -0.019607843137254943, ''MMH.NZ'': 0.027732463295269127, ''MPG.NZ'': -0.011235955056179803, ''MWE.NZ'': 0.04081632653061229, ''NTL.NZ'': 0.0, ''NWF.NZ'': -0.017857142857142905, ''NZK.NZ'': -0.006896551724137945, ''NZM.NZ'': 0.010638297872340496, ''NZO.NZ'': 0.0, ''NZR.NZ'': 0.0, ''NZX.NZ'': -0.021505376344086002, ''OCA.NZ'': 0.006666666666666599, ''PCT.NZ'': -0.002941176470588225, ''PEB.NZ'': -0.03597122302158262, ''PFI.NZ'': -0.009900990099009799, ''PGW.NZ'': -0.010752688172043001, ''PHL.NZ'': 0.0, ''POT.NZ'': 0.03216783216783203, ''PPH.NZ'': 0.005714285714285783, ''PX1.NZ'': -0.05263157894736825, ''PYS.NZ'': -0.015384615384615441, ''QEX.NZ'': 0.0, ''RAK.NZ'': -0.020618556701030966, ''RBD.NZ'': 0.024983984625240208, ''RYM.NZ'': 0.011235955056179803, ''SAN.NZ'': -0.008620689655172376, ''SCL.NZ'': 0.030303030303030276, ''SCT.NZ'': 0.03191489361702149, ''SCY.NZ'': 0.0, ''SDL.NZ'': 0.03333333333333344, ''SEK.NZ'': 0.003824091778202643, ''SKC.NZ'': 0.01834862385321112, ''SKL.NZ'': 0.03612167300380231, ''SKO.NZ'': -0.0024937655860348684, ''SKT.NZ'': 0.01840490797546024, ''SML.NZ'': -0.015243902439024293, ''SPG.NZ'': -0.015151515151515138, ''SPK.NZ'': 0.0, ''SPN.NZ'': 0.008235294117647118, ''SPY.NZ'': 0.01333333333333342, ''STU.NZ'': 0.0, ''SUM.NZ'': 0.009210526315789558, ''TEM.NZ'': 0.02556818181818188, ''TGG.NZ'': 0.0, ''THL.NZ'': 0.04721030042918439, ''TLL.NZ'': 0.0, ''TLT.NZ'': 0.0, ''TPW.NZ'': -0.018633540372670843, ''TRA.NZ'': 0.0, ''TRU.NZ'': 0.03333333333333344, ''TWR.NZ'': -0.0071428571428571175, ''VCT.NZ'': -0.025821596244131273, ''VGL.NZ'': 0.08370044052863435, ''VHP.NZ'': -0.006069802731411222, ''VTL.NZ'': 0.013513513513513598, ''WBC.NZ'': -0.0007382798080471931, ''WDT.NZ'': -0.03149606299212604, ''WHS.NZ'': -0.0026385224274406704, ''ZEL.NZ'': 0.0}')
    ,@(4, ' 29/08/2021', -0.004358792098968958, -0.0009567564055181754, 0.001931470052109132, 0.004986301369863014, '{''ACE.NZ'': 0.0, ''AFC.NZ'': 0.0, ''AFI.NZ'': -0.0034207525655644, ''AFT.NZ'': -0.011494252873563204, ''AGL.NZ'': 0.024844720496894235, ''AIA.NZ'': -0.0070224719101122934, ''AIR.NZ'': -0.00666666666666671, ''ALF.NZ'': -0.016129032258064502, ''AMP.NZ'': 0.0, ''ANZ.NZ'': -0.003712453594330034, ''APL.NZ'': -0.014925373134328401, ''ARB.NZ'': -0.032258064516129115, ''ARG.NZ'': -0.006134969325153228, ''ARV.NZ'': 0.004878048780487809, ''ATM.NZ'': -0.039432176656151396, ''BFG.NZ'': -0.012658227848101222, ''BGI.NZ'': 0.0, ''BGP.NZ'': -0.0015384615384614886, ''BIT.NZ'': -0.029787234042553346, ''BLT.NZ'': 0.0, ''BRM.NZ'': 0.0, ''CAV.NZ'': 0.01538461538461533, ''CBD.NZ'': 0.012195121951219523, ''CDI.NZ'': 0.01754385964912286, ''CEN.NZ'': -0.008663366336633671, ''CGF.NZ'': 0.026315789473684292, ''CMO.NZ'': 0.0, ''CNU.NZ'': 0.013138686131386912, ''CRP.NZ'': 0.0, ''DGL.NZ'': 0.015602836879432758, ''EBO.NZ'': 0.0005712653527565337, ''ERD.NZ'': 0.0031847133757960666, ''EVO.NZ'': 0.0, ''FBU.NZ'': -0.02785145888594165, ''FCT.NZ'': -0.02506963788300831, ''FPH.NZ'': 0.008442776735459567, ''FRE.NZ'': -0.011032308904649235, ''FSF.NZ'': 0.0, ''FWL.NZ'': 0.013605442176870763, ''GEN.NZ'': 0.0, ''GEO.NZ'': -0.025423728813559254, ''GFL.NZ'': 0.01388888888888884, ''GMT.NZ'': 0.011450381679389166, ''GNE.NZ'': 0.0, ''GSH.NZ'': -0.014084507042253391, ''GTK.NZ'': -0.00990099009900991, ''GWC.NZ'': 0.0, ''GXH.NZ'': 0.017699115044247815, ''HFL.NZ'': 0.0, ''HGH.NZ'': 0.008928571428571175, ''HLG.NZ'': -0.0098314606741573, ''IFT.NZ'': 0.0054421768707482165, ''IKE.NZ'': 0.0, ''IPL.NZ'': -0.009803921568627416, ''JLG.NZ'': -0.010869565217391353, ''JPG.NZ'': -0.004319654427645703, ''KFL.NZ'': 0.0, ''KMD.NZ'': -0.014705882352941235, ''KPG.NZ'': 0.008658008658008587, ''MCK.NZ'': -0.012295081967213073, ''MCY.NZ'': 0.0, ''MEL.NZ'': 0.009960159362550014, ''MFT.NZ'': 0.005440104450005423, ''MHJ.NZ'': 0.0, ''MHM.NZ'': 0.0, ''MLN.NZ'': -0.019607843137254943, ''MMH.NZ'': 0.027732463295269127, ''MPG.NZ'': -0.011235955056179803, ''MWE.NZ'': 0.04081632653061229, ''NTL.NZ'': 0.0, ''NWF.NZ'': -0.017857142857142905, ''NZK.NZ'': -0.006896551724137945, ''NZM.NZ'': 0.010638297872340496, ''NZO.NZ'': 0.0, ''NZR.NZ'': 0.0, ''NZX.NZ'': -0.021505376344086002, ''OCA.NZ'': 0.006666666666666599, ''PCT.NZ'': -0.002941176470588225, ''PEB.NZ'': -0.03597122302158262, ''PFI.NZ'': -0.009900990099009799, ''PGW.NZ'': -0.010752688172043001, ''PHL.NZ'': 0.0, ''POT.NZ'': 0.03216783216783203, ''PPH.NZ'': 0.005714285714285783, ''PX1.NZ'': -0.05263157894736825, ''PYS.NZ'': -0.015384615384615441, ''QEX.NZ'': 0.0, ''RAK.NZ'': -0.020618556701030966, ''RBD.NZ'': 0.024983984625240208, ''RYM.NZ'': 0.011235955056179803, ''SAN.NZ'': -0.008620689655172376, ''SCL.NZ'': 0.030303030303030276, ''SCT.NZ'': 0.03191489361702149, ''SCY.NZ'': 0.0, ''SDL.NZ'': 0.03333333333333344, ''SEK.NZ'': 0.003824091778202643, ''SKC.NZ'': 0.01834862385321112, ''SKL.NZ'': 0.03612167300380231, ''SKO.NZ'': -0.0024937655860348684, ''SKT.NZ'': 0.01840490797546024, ''SML.NZ'': -0.015243902439024293, ''SPG.NZ'': -0.015151515151515138, ''SPK.NZ'': 0.0, ''SPN.NZ'': 0.008235294117647118, ''SPY.NZ'': 0.01333333333333342, ''STU.NZ'': 0.0, ''SUM.NZ'': 0.009210526315789558, ''TEM.NZ'': 0.02556818181818188, ''TGG.NZ'': 0.0, ''THL.NZ'': 0.04721030042918439, ''TLL.NZ'': 0.0, ''TLT.NZ'': 0.0, ''TPW.NZ'': -0.018633540372670843, ''TRA.NZ'': 0.0, ''TRU.NZ'': 0.03333333333333344, ''TWR.NZ'': -0.0071428571428571175, ''VCT.NZ'': -0.025821596244131273, ''VGL.NZ'': 0.08370044052863435, ''VHP.NZ'': -0.006069802731411222, ''VTL.NZ'': 0.013513513513513598, ''WBC.NZ'': -0.0007382798080471931, ''WDT.NZ'': -0.03149606299212604, ''WHS.NZ'': -0.0026385224274406704, ''ZEL.NZ'': 0.0}')
)

$startRow = 3

# Copy the style of the existing A2 cell (bold/centered header-like style)
# so the new index cells in column A (A3:A6) match it, same as the template row.
$styleSource = $ws.Cells.Item(2, 1)
$styleSource.Copy()

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $data = $rowsData[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]

    $ws.Cells.Item($r, 1).PasteSpecial(-4122) # xlPasteFormats, matches A2 style (s="1")
}

Write-Host "Added rows 3-6"